$p = $ppt.ActivePresentation
Write-Output ($p.SlideMaster -eq $null)
